$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("WV50 FILTER", "Traza", "2024-06-10", "10:35:35", "Mañana", "10:35:36", "0:00:01", "-0.00 minutos"),
    @("WV50 FILTER", "Robot no coge PCB", "2024-06-10", "10:35:38", "Mañana", "10:35:40", "0:00:02", "0.01 minutos"),
    @("WV50 FILTER", "NOK Soldadura Plástico", "2024-06-10", "10:36:59", "Mañana", "10:37:00", "0:00:01", "0.23 minutos"),
    @("SPL", "No detecta marcas Power", "2024-06-10", "10:48:07", "Mañana", "10:48:09", "0:00:02", "-0.00 minutos"),
    @("SPL", "Soldadura defectuosa", "2024-06-10", "10:48:10", "Mañana", "10:48:12", "0:00:02", "0.01 minutos"),
    @("SPL", "No detecta marcas Power", "2024-06-10", "10:48:43", "Mañana", "10:48:46", "0:00:03", "0.10 minutos"),
    @("SPL", "Error en sensor de salida", "2024-06-10", "10:48:51", "Mañana", "10:48:52", "0:00:01", "0.09 minutos")
)

$startRow = 159
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $value = $rowData[$c]
        if ($c -eq 2) {
            # Column C holds a date-formatted string (e.g. "2024-06-10") that
            # must remain literal text rather than be auto-converted to a
            # date serial number, so force it with a text quote-prefix.
            $ws.Cells.Item($r, $c + 1).Value = "'" + $value
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $value
        }
    }
}
